# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-name suffixes to "_FV2310" / "_FV2404"
# - Turn the data range into a real Excel Table (ListObject)
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row labels: *_old -> *_FV2310, *_new -> *_FV2404
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($header -ne $null) {
        $renamed = $header -replace "_old$", "_FV2310"
        $renamed = $renamed -replace "_new$", "_FV2404"
        if ($renamed -ne $header) {
            $cell.Value = $renamed
        }
    }
}

# 2. Convert the used range A1:U73 into an Excel Table
$dataRange = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
